$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105, pushing existing rows 105-119 down to 106-120.
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new data record.
$ws.Range("A105").Value = 9
$ws.Range("B105").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C105").Value = "Metropolitana"
$ws.Range("D105").Value = 44511
$ws.Range("E105").Value = 13
$ws.Range("F105").Value = "Fruta"
$ws.Range("G105").Value = 100101
$ws.Range("H105").Value = "Berries"
$ws.Range("I105").Value = 100101001
$ws.Range("J105").Value = "Arándano (blue)"
$ws.Range("K105").Value = "Sin especificar"
$ws.Range("L105").Value = "Primera"
$ws.Range("M105").Value = 500
$ws.Range("N105").Value = 5000
$ws.Range("O105").Value = 5000
$ws.Range("P105").Value = 5000
$ws.Range("Q105").Value = "$/bandeja 2 kilos"
$ws.Range("R105").Value = "Región Metropolitana"
$ws.Range("S105").Value = 2500
$ws.Range("T105").Value = 2
